$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 439
$ws.Range("F7").Value = 1289
$ws.Range("F8").Value = 485
$ws.Range("F10").Value = 255
$ws.Range("F11").Value = 166
$ws.Range("F12").Value = 203
$ws.Range("F13").Value = 1084
$ws.Range("F14").Value = 6
$ws.Range("F17").Value = 51
$ws.Range("F18").Value = 229
$ws.Range("F19").Value = 1611
$ws.Range("F20").Value = 593
$ws.Range("F21").Value = 257
$ws.Range("F22").Value = 128
$ws.Range("F23").Value = 1048
$ws.Range("F24").Value = 388
$ws.Range("F26").Value = 906
$ws.Range("F27").Value = 1186
$ws.Range("F29").Value = 1905
$ws.Range("F30").Value = 2775
$ws.Range("F31").Value = 1549
$ws.Range("F33").Value = 88
$ws.Range("F34").Value = 572
$ws.Range("F35").Value = 846
$ws.Range("F36").Value = 1583
$ws.Range("F37").Value = 865
$ws.Range("F38").Value = 1631
$ws.Range("F39").Value = 187
$ws.Range("F41").Value = 818
$ws.Range("F42").Value = 13
$ws.Range("F43").Value = 764
$ws.Range("F44").Value = 748
$ws.Range("F45").Value = 954
$ws.Range("F46").Value = 410
$ws.Range("F47").Value = 3289

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 10
$ws.Range("F15").Value = 763

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 439
$ws.Range("F10").Value = 1289
$ws.Range("F11").Value = 485
$ws.Range("F13").Value = 255
$ws.Range("F14").Value = 166
$ws.Range("F15").Value = 203
$ws.Range("F16").Value = 1084
$ws.Range("F19").Value = 51
$ws.Range("F20").Value = 229
$ws.Range("F21").Value = 1611
$ws.Range("F22").Value = 593
$ws.Range("F23").Value = 257
$ws.Range("F24").Value = 388
$ws.Range("F25").Value = 10
$ws.Range("F28").Value = 1186
$ws.Range("F29").Value = 2775
$ws.Range("F31").Value = 1549
$ws.Range("F33").Value = 763
$ws.Range("F35").Value = 572
$ws.Range("F36").Value = 846
$ws.Range("F37").Value = 1583
$ws.Range("F39").Value = 865
$ws.Range("F40").Value = 1631
$ws.Range("F41").Value = 818
$ws.Range("F42").Value = 764
$ws.Range("F43").Value = 748
$ws.Range("F44").Value = 954
$ws.Range("F45").Value = 410
$ws.Range("F48").Value = 3289
